$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed VOL/B0/CET_100/CET_300/SP2_PERC/SP/SP2/SP3 for existing rows (2-26),
# plus five new CRN samples appended as rows 27-31 (crn_00025..crn_00029).
# Data rows: @(row, B, C, D, E, F, G, H, I)
$data = @(
    @(2, 871.0112360000001, 257.805191, 1.107013, 7.83053, 50, 2, 64, 62),
    @(3, 873.174569, 231.377293, 0.029002, 4.592623, 53.125, 1, 68, 59),
    @(4, 862.094684, 250.665492, -2.369145, 0.680936, 45.3125, 1, 58, 69),
    @(5, 867.26811, 220.789683, 0.173808, 7.256983, 54.6875, 3, 70, 55),
    @(6, 842.596559, 234.21166, 5.926872, 34.730584, 54.6875, 2, 70, 56),
    @(7, 904.960629, 254.428844, 2.863191, 15.61796, 50, 3, 64, 61),
    @(8, 844.846649, 269.202894, 0.099444, 4.735319, 51.5625, 1, 66, 61),
    @(9, 864.188796, 247.48918, 0.032859, 2.498639, 45.3125, 3, 58, 67),
    @(10, 862.121858, 256.141665, 0.6160870000000001, 4.858, 50, 2, 64, 62),
    @(11, 856.076057, 241.177543, 2.153683, 11.602145, 54.6875, 1, 70, 57),
    @(12, 853.887066, 283.605365, 0.310456, 5.5937, 46.875, 2, 60, 66),
    @(13, 873.448664, 267.758317, -0.011156, 5.760397, 46.875, 1, 60, 67),
    @(14, 877.59447, 256.374091, -0.478993, 2.114145, 53.125, 3, 68, 57),
    @(15, 860.248153, 251.468366, 0.323785, 3.561581, 54.6875, 2, 70, 56),
    @(16, 861.279713, 272.54073, 0.201327, 4.48071, 51.5625, 2, 66, 60),
    @(17, 883.299083, 235.417525, 2.339854, 30.465007, 50, 3, 64, 61),
    @(18, 856.319508, 278.482293, -0.509336, 2.914026, 45.3125, 0, 58, 70),
    @(19, 846.026864, 263.035889, 0.065183, 5.93822, 50, 1, 64, 63),
    @(20, 881.367229, 213.648102, -0.569958, 1.340945, 54.6875, 3, 70, 55),
    @(21, 892.472232, 247.529616, -0.307756, 2.603753, 50, 2, 64, 62),
    @(22, 864.4337839999999, 247.761817, 0.040677, 2.494565, 51.5625, 2, 66, 60),
    @(23, 877.1545589999999, 273.594081, 0.056319, 3.804761, 50, 1, 64, 63),
    @(24, 869.315935, 242.165303, 0.200015, 8.358812, 50, 2, 64, 62),
    @(25, 860.7095389999999, 277.149964, 0.158575, 3.952814, 46.875, 2, 60, 66),
    @(26, 870.8744369999999, 238.912853, 0.355134, 13.098271, 50, 3, 64, 61),
    @(27, 867.404531, 295.248265, -0.919385, 5.539313, 46.875, 0, 60, 68),
    @(28, 882.305206, 218.405823, 1.339119, 9.654825000000001, 54.6875, 3, 70, 55),
    @(29, 871.8049549999999, 239.230092, -0.289547, 2.65072, 45.3125, 2, 58, 68),
    @(30, 859.4885839999999, 233.068891, 1.035501, 8.482517, 50, 2, 64, 62),
    @(31, 867.540967, 277.950326, 0.349355, 8.73907, 48.4375, 1, 62, 65)
)

# New rows need the same formatting (borders/fill) as the existing data rows.
$ws.Range("A2:I2").Copy()
$ws.Range("A27:I31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

foreach ($row in $data) {
    $r = $row[0]
    if ($r -ge 27) {
        $crnIndex = $r - 2
        $crnText = "crn_{0:D5}" -f $crnIndex
        $ws.Range("A$r").Value = $crnText
    }
    $arr = New-Object 'object[,]' 1,8
    $arr[0,0] = $row[1]
    $arr[0,1] = $row[2]
    $arr[0,2] = $row[3]
    $arr[0,3] = $row[4]
    $arr[0,4] = $row[5]
    $arr[0,5] = $row[6]
    $arr[0,6] = $row[7]
    $arr[0,7] = $row[8]
    $ws.Range("B${r}:I${r}").Value = $arr
}

Write-Output "Done updating sheet1 rows 2-31"